$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.223.16'
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = '  -6.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.213.88'
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = '  -6.39%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.94'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  -6.75%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.93'
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = '  -5.37%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.552'
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = '  -8.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.45'
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = '  -7.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '57.81'
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = '  -3.22%  '

$ws.Range("E13").Value = '  -4.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.70'
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = '  -7.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.542.75'
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = '  -6.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.75'
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = '  -9.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.839'
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = '  -9.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.211.87'
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = '  -6.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.152.73'
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = '  -6.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = '  -8.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.30'
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = '  -6.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.08'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  -7.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.51'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  -8.86%  '

$ws.Range("E24").Value = '  +6.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("E26").Value = '  -5.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.42'
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = '  -3.01%  '

$ws.Range("E28").Value = '  -4.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.76'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  -7.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.78'
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = '  -2.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.43'
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = '  -8.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.119'
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = '  -8.41%  '

$ws.Range("E33").Value = '  -7.75%  '

$ws.Range("E34").Value = '  -6.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.16'
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = '  -4.62%  '

$ws.Range("E36").Value = '  -9.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.89'
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = '  +2.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.73'
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = '  +15.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0277'
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = '  -1.85%  '

$ws.Range("E40").Value = '  -5.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.82'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  -11.98%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.00'
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = '  -2.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.98'
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = '  -10.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.196'
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.61'
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = '  -5.00%  '

$ws.Range("E46").Value = '  -6.38%  '

$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.63'
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  +10.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.44'
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = '  +1.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = '  -6.37%  '

$ws.Range("E51").Value = '  -5.57%  '
